$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 3 so the "mariusz" entry is split into two rows
# (row 2 and row 3), pushing the former rows 3 (user123) and 4 (user1234)
# down to rows 4 and 5.
$ws.Rows.Item(3).Insert()

# Row 2: mariusz / testDocumentMariusz1.docx
$ws.Range("A2").Value = "mariusz"
$ws.Range("B2").Value = "/documentsList/documents/testDocumentMariusz1.docx"

# Row 3: mariusz / testDocumentMariusz2.docx
$ws.Range("A3").Value = "mariusz"
$ws.Range("B3").Value = "/documentsList/documents/testDocumentMariusz2.docx"

# Row 4: user123 / testDocumentUser123.docx (was row 3)
$ws.Range("A4").Value = "user123"
$ws.Range("B4").Value = "/documentsList/documents/testDocumentUser123.docx"

# Row 5: user1234 / testDocumentUser1234.pdf (was row 4)
$ws.Range("A5").Value = "user1234"
$ws.Range("B5").Value = "/documentsList/documents/testDocumentUser1234.pdf"

# Update the selection to match the diff (active cell on B4)
$ws.Range("B4").Select()
